# Update the NATMI LR-pairs worksheet with the new TPM-derived results.
# The new data keeps only 3 result rows (previously rows 5-7, sender = MuSCs)
# with refreshed numeric values, and the old rows 2-4 (sender = ECs) are removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-obsolete rows 5 through 7 first (row 2-4 values are overwritten below).
$ws.Rows("5:7").Delete()

# Row 2: MuSCs -> Osm -> Il6st -> ECs
$ws.Range("A2").Value = "MuSCs"
$ws.Range("B2").Value = "Osm"
$ws.Range("C2").Value = "Il6st"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.03867233333333333
$ws.Range("H2").Value = 0.116017
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 20.854426
$ws.Range("N2").Value = 62.563278
$ws.Range("O2").Value = 0.1507164072139519
$ws.Range("P2").Value = 0.1507164072139519
$ws.Range("Q2").Value = 0.8064893137473332
$ws.Range("R2").Value = 7.258403823725999
$ws.Range("S2").Value = 0.1507164072139519
$ws.Range("T2").Value = 0.1507164072139519

# Row 3: MuSCs -> Osm -> Il6st -> FAPs
$ws.Range("A3").Value = "MuSCs"
$ws.Range("B3").Value = "Osm"
$ws.Range("C3").Value = "Il6st"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.03867233333333333
$ws.Range("H3").Value = 0.116017
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 94.96115633333334
$ws.Range("N3").Value = 284.883469
$ws.Range("O3").Value = 0.6862909728343718
$ws.Range("P3").Value = 0.6862909728343718
$ws.Range("Q3").Value = 3.672369491441444
$ws.Range("R3").Value = 33.051325422973
$ws.Range("S3").Value = 0.6862909728343718
$ws.Range("T3").Value = 0.6862909728343718

# Row 4: MuSCs -> Osm -> Il6st -> MuSCs
$ws.Range("A4").Value = "MuSCs"
$ws.Range("B4").Value = "Osm"
$ws.Range("C4").Value = "Il6st"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.03867233333333333
$ws.Range("H4").Value = 0.116017
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 22.553069
$ws.Range("N4").Value = 67.65920700000001
$ws.Range("O4").Value = 0.1629926199516763
$ws.Range("P4").Value = 0.1629926199516763
$ws.Range("Q4").Value = 0.8721798020576668
$ws.Range("R4").Value = 7.849618218519001
$ws.Range("S4").Value = 0.1629926199516763
$ws.Range("T4").Value = 0.1629926199516763
